# TestSubjectDataExportTemplate.xlsx edit
#
# The only real content change in this revision is a text fix in the
# "Measurements" table header: "Maximum contraction [Nm] / Fatigue [%]"
# becomes "Maximum contraction [N] / Fatigue [%]" (the unit changes from
# Nm to N). Everything else in the underlying OOXML diff (shared-string /
# dxf re-ordering, GUIDs, rupBuild, absPath, window size, cursor position)
# is Excel's own resave churn and carries no semantic meaning.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header cell of the "Maximum contraction" column (column E,
# row 6 of the MeasurementsTable) - this also renames the table column
# since it's the table's header row.
$ws.Range("E6").Value2 = "Maximum contraction [N] `n/ Fatigue [%] "

# Move the active selection the way the author left it when saving.
$ws.Range("F11").Select()
